$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data of rows 2<->4 and 3<->5 for columns D, M, N, O, P, R, S
# (A, B, C, E, F, G, H, I, J, K, L, Q, T are identical across these rows, so no change needed there)

$row2 = @{ D = 44320; M = 50;  N = 18000; O = 20000; P = 18800; R = "Provincia de Limarí";  S = 1044 }
$row3 = @{ D = 44719; M = 50;  N = 20000; O = 21000; P = 20400; R = "Provincia de Limarí";  S = 1133 }
$row4 = @{ D = 44362; M = 100; N = 19000; O = 20000; P = 19500; R = "Provincia de Curicó";  S = 1083 }
$row5 = @{ D = 45084; M = 100; N = 17000; O = 18000; P = 17500; R = "Región de O'Higgins"; S = 972 }

function Set-RowValues($rowNum, $values) {
    $ws.Range("D$rowNum").Value = $values.D
    $ws.Range("M$rowNum").Value = $values.M
    $ws.Range("N$rowNum").Value = $values.N
    $ws.Range("O$rowNum").Value = $values.O
    $ws.Range("P$rowNum").Value = $values.P
    $ws.Range("R$rowNum").Value = $values.R
    $ws.Range("S$rowNum").Value = $values.S
}

Set-RowValues 2 $row2
Set-RowValues 3 $row3
Set-RowValues 4 $row4
Set-RowValues 5 $row5
